$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 351

# New row 3
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 287

# New row 4
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 234

# New row 5
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 130

# Copy the style from A2 (existing bold/border/centered style) to the new A3:A5 cells
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122) | Out-Null
